$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WORK")

$ws.Range("B7").Value = 59000000.0
$ws.Range("C7").Value = 63000000.0
$ws.Range("D7").Value = 24121000.0
$ws.Range("E7").Value = 46177000.0
$ws.Range("F7").Value = 17879000.0
